$d = $word.ActiveDocument

$pairs = @(
  @("799÷3=", "326÷8="),
  @("579÷8=", "617÷3="),
  @("870÷8=", "333÷9="),
  @("455÷6=", "117÷7="),
  @("919÷7=", "400÷8="),
  @("425÷2=", "120÷3="),
  @("294÷8=", "987÷9="),
  @("741÷6=", "491÷9="),
  @("448÷7=", "114÷4="),
  @("284÷8=", "764÷8="),
  @("577÷6=", "337÷5="),
  @("135÷4=", "592÷8="),
  @("397÷3=", "123÷2="),
  @("467÷7=", "710÷5="),
  @("647÷7=", "681÷5="),
  @("375÷9=", "893÷3="),
  @("914÷4=", "413÷3="),
  @("585÷2=", "661÷4="),
  @("133÷5=", "989÷5="),
  @("207÷9=", "282÷2="),
  @("682÷7=", "808÷6="),
  @("321÷5=", "991÷6="),
  @("512÷2=", "775÷5="),
  @("703÷4=", "352÷6="),
  @("882÷6=", "310÷5=")
)

foreach ($pair in $pairs) {
    $old = $pair[0]
    $new = $pair[1]
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false,
                             $true, 1, $false, $new, 2)
}
